$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New time-log entries (row 55 intentionally left blank, matching the
# existing pattern of skipped rows in this sheet).
$ws.Range("A54").Value = 45247
$ws.Range("A54").NumberFormat = "d-mmm"
$ws.Range("B54").Value = 3
$ws.Range("C54").Value = "started implementing the integration between the front end and back, starting with printing the travel advisories on the panels on the UI section< but its still very buggy and it doesn’t work"

$ws.Range("A56").Value = 45248
$ws.Range("A56").NumberFormat = "d-mmm"
$ws.Range("B56").Value = 3
$ws.Range("C56").Value = "debugging integrations between frontend and backend, thinking either sending everything to database then fetching, just printing it all out front. "

$ws.Range("A57").Value = 45249
$ws.Range("A57").NumberFormat = "d-mmm"
$ws.Range("B57").Value = 4
$ws.Range("C57").Value = "budget panel frontend and backend integrated however, the code is still buggy, continuing to work on it"

$ws.Range("C57").Select()
